$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D191").Value = 44868
$ws.Range("K191").Value = "Lane Late"
$ws.Range("L191").Value = "Primera"
$ws.Range("M191").Value = 160
$ws.Range("N191").Value = 8000
$ws.Range("O191").Value = 9000
$ws.Range("P191").Value = 8500
$ws.Range("Q191").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R191").Value = "Región de O'Higgins"
$ws.Range("S191").Value = 567
$ws.Range("T191").Value = 15

$ws.Range("D192").Value = 44477
$ws.Range("K192").Value = "Lane Late"
$ws.Range("L192").Value = "Primera"
$ws.Range("M192").Value = 240
$ws.Range("N192").Value = 6000
$ws.Range("O192").Value = 6500
$ws.Range("P192").Value = 6250
$ws.Range("Q192").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R192").Value = "Región de O'Higgins"
$ws.Range("S192").Value = 417
$ws.Range("T192").Value = 15

$ws.Range("D193").Value = 44477
$ws.Range("K193").Value = "Lane Late"
$ws.Range("L193").Value = "Segunda"
$ws.Range("M193").Value = 240
$ws.Range("N193").Value = 5000
$ws.Range("O193").Value = 5500
$ws.Range("P193").Value = 5250
$ws.Range("Q193").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R193").Value = "Región de O'Higgins"
$ws.Range("S193").Value = 350
$ws.Range("T193").Value = 15

$ws.Range("D194").Value = 44477
$ws.Range("K194").Value = "Navel Late"
$ws.Range("L194").Value = "Primera"
$ws.Range("M194").Value = 240
$ws.Range("N194").Value = 6000
$ws.Range("O194").Value = 6500
$ws.Range("P194").Value = 6250
$ws.Range("Q194").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R194").Value = "Región de O'Higgins"
$ws.Range("S194").Value = 417
$ws.Range("T194").Value = 15

$ws.Range("D195").Value = 44477
$ws.Range("K195").Value = "Navel Late"
$ws.Range("L195").Value = "Segunda"
$ws.Range("M195").Value = 240
$ws.Range("N195").Value = 5000
$ws.Range("O195").Value = 5500
$ws.Range("P195").Value = 5250
$ws.Range("Q195").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R195").Value = "Región de O'Higgins"
$ws.Range("S195").Value = 350
$ws.Range("T195").Value = 15

$ws.Range("D196").Value = 44426
$ws.Range("K196").Value = "Fukumoto"
$ws.Range("L196").Value = "Primera"
$ws.Range("M196").Value = 160
$ws.Range("N196").Value = 6000
$ws.Range("O196").Value = 6200
$ws.Range("P196").Value = 6100
$ws.Range("Q196").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R196").Value = "Región de O'Higgins"
$ws.Range("S196").Value = 407
$ws.Range("T196").Value = 15

$ws.Range("D197").Value = 44426
$ws.Range("K197").Value = "Fukumoto"
$ws.Range("L197").Value = "Segunda"
$ws.Range("M197").Value = 70
$ws.Range("N197").Value = 5000
$ws.Range("O197").Value = 5500
$ws.Range("P197").Value = 5214
$ws.Range("Q197").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R197").Value = "Región de O'Higgins"
$ws.Range("S197").Value = 348
$ws.Range("T197").Value = 15

$ws.Range("D198").Value = 44426
$ws.Range("K198").Value = "Navel Late"
$ws.Range("L198").Value = "Primera"
$ws.Range("M198").Value = 160
$ws.Range("N198").Value = 5000
$ws.Range("O198").Value = 5500
$ws.Range("P198").Value = 5250
$ws.Range("Q198").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R198").Value = "Región de O'Higgins"
$ws.Range("S198").Value = 350
$ws.Range("T198").Value = 15

$ws.Range("D199").Value = 44426
$ws.Range("K199").Value = "Navel Late"
$ws.Range("L199").Value = "Segunda"
$ws.Range("M199").Value = 80
$ws.Range("N199").Value = 4500
$ws.Range("O199").Value = 4800
$ws.Range("P199").Value = 4650
$ws.Range("Q199").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R199").Value = "Región de O'Higgins"
$ws.Range("S199").Value = 310
$ws.Range("T199").Value = 15

$ws.Range("D200").Value = 44658
$ws.Range("K200").Value = "Valencia"
$ws.Range("L200").Value = "Primera"
$ws.Range("M200").Value = 160
$ws.Range("N200").Value = 9500
$ws.Range("O200").Value = 10000
$ws.Range("P200").Value = 9750
$ws.Range("Q200").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R200").Value = "Región de O'Higgins"
$ws.Range("S200").Value = 650
$ws.Range("T200").Value = 15

$ws.Range("D201").Value = 44508
$ws.Range("K201").Value = "Lane Late"
$ws.Range("L201").Value = "Primera"
$ws.Range("M201").Value = 160
$ws.Range("N201").Value = 8000
$ws.Range("O201").Value = 8500
$ws.Range("P201").Value = 8250
$ws.Range("Q201").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R201").Value = "Región de O'Higgins"
$ws.Range("S201").Value = 550
$ws.Range("T201").Value = 15

$ws.Range("D202").Value = 44508
$ws.Range("K202").Value = "Lane Late"
$ws.Range("L202").Value = "Segunda"
$ws.Range("M202").Value = 120
$ws.Range("N202").Value = 7000
$ws.Range("O202").Value = 7500
$ws.Range("P202").Value = 7250
$ws.Range("Q202").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R202").Value = "Región de O'Higgins"
$ws.Range("S202").Value = 483
$ws.Range("T202").Value = 15

$ws.Range("D203").Value = 44629
$ws.Range("K203").Value = "Valencia"
$ws.Range("L203").Value = "Primera"
$ws.Range("M203").Value = 80
$ws.Range("N203").Value = 11000
$ws.Range("O203").Value = 12000
$ws.Range("P203").Value = 11500
$ws.Range("Q203").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R203").Value = "Región de O'Higgins"
$ws.Range("S203").Value = 767
$ws.Range("T203").Value = 15

$ws.Range("D204").Value = 44629
$ws.Range("K204").Value = "Valencia"
$ws.Range("L204").Value = "Segunda"
$ws.Range("M204").Value = 60
$ws.Range("N204").Value = 9000
$ws.Range("O204").Value = 9000
$ws.Range("P204").Value = 9000
$ws.Range("Q204").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R204").Value = "Región de O'Higgins"
$ws.Range("S204").Value = 600
$ws.Range("T204").Value = 15

$ws.Range("D205").Value = 44762
$ws.Range("K205").Value = "Fukumoto"
$ws.Range("L205").Value = "Primera"
$ws.Range("M205").Value = 120
$ws.Range("N205").Value = 6500
$ws.Range("O205").Value = 7000
$ws.Range("P205").Value = 6750
$ws.Range("Q205").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R205").Value = "Región de O'Higgins"
$ws.Range("S205").Value = 450
$ws.Range("T205").Value = 15

$ws.Range("D206").Value = 44762
$ws.Range("K206").Value = "Fukumoto"
$ws.Range("L206").Value = "Segunda"
$ws.Range("M206").Value = 60
$ws.Range("N206").Value = 6000
$ws.Range("O206").Value = 6000
$ws.Range("P206").Value = 6000
$ws.Range("Q206").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R206").Value = "Región de O'Higgins"
$ws.Range("S206").Value = 400
$ws.Range("T206").Value = 15

$ws.Range("D207").Value = 44754
$ws.Range("K207").Value = "Fukumoto"
$ws.Range("L207").Value = "Primera"
$ws.Range("M207").Value = 100
$ws.Range("N207").Value = 7500
$ws.Range("O207").Value = 8000
$ws.Range("P207").Value = 7750
$ws.Range("Q207").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R207").Value = "Región de O'Higgins"
$ws.Range("S207").Value = 517
$ws.Range("T207").Value = 15

$ws.Range("D208").Value = 44754
$ws.Range("K208").Value = "Fukumoto"
$ws.Range("L208").Value = "Segunda"
$ws.Range("M208").Value = 50
$ws.Range("N208").Value = 6000
$ws.Range("O208").Value = 6000
$ws.Range("P208").Value = 6000
$ws.Range("Q208").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R208").Value = "Región de O'Higgins"
$ws.Range("S208").Value = 400
$ws.Range("T208").Value = 15

$ws.Range("D209").Value = 44673
$ws.Range("K209").Value = "Valencia"
$ws.Range("L209").Value = "Primera"
$ws.Range("M209").Value = 120
$ws.Range("N209").Value = 10000
$ws.Range("O209").Value = 11000
$ws.Range("P209").Value = 10500
$ws.Range("Q209").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R209").Value = "Región de O'Higgins"
$ws.Range("S209").Value = 700
$ws.Range("T209").Value = 15

$ws.Range("D210").Value = 44196
$ws.Range("K210").Value = "Valencia"
$ws.Range("L210").Value = "Primera"
$ws.Range("M210").Value = 150
$ws.Range("N210").Value = 16000
$ws.Range("O210").Value = 17000
$ws.Range("P210").Value = 16467
$ws.Range("Q210").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R210").Value = "Región de O'Higgins"
$ws.Range("S210").Value = 1098
$ws.Range("T210").Value = 15

$ws.Range("D211").Value = 44677
$ws.Range("K211").Value = "Valencia"
$ws.Range("L211").Value = "Primera"
$ws.Range("M211").Value = 160
$ws.Range("N211").Value = 10000
$ws.Range("O211").Value = 11000
$ws.Range("P211").Value = 10500
$ws.Range("Q211").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R211").Value = "Región de O'Higgins"
$ws.Range("S211").Value = 700
$ws.Range("T211").Value = 15

$ws.Range("D212").Value = 44677
$ws.Range("K212").Value = "Valencia"
$ws.Range("L212").Value = "Segunda"
$ws.Range("M212").Value = 80
$ws.Range("N212").Value = 9000
$ws.Range("O212").Value = 9000
$ws.Range("P212").Value = 9000
$ws.Range("Q212").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R212").Value = "Región de O'Higgins"
$ws.Range("S212").Value = 600
$ws.Range("T212").Value = 15

$ws.Range("D213").Value = 44509
$ws.Range("K213").Value = "Lane Late"
$ws.Range("L213").Value = "Primera"
$ws.Range("M213").Value = 160
$ws.Range("N213").Value = 8000
$ws.Range("O213").Value = 8500
$ws.Range("P213").Value = 8250
$ws.Range("Q213").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R213").Value = "Región de O'Higgins"
$ws.Range("S213").Value = 550
$ws.Range("T213").Value = 15

$ws.Range("D214").Value = 44509
$ws.Range("K214").Value = "Lane Late"
$ws.Range("L214").Value = "Segunda"
$ws.Range("M214").Value = 120
$ws.Range("N214").Value = 7000
$ws.Range("O214").Value = 7500
$ws.Range("P214").Value = 7250
$ws.Range("Q214").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R214").Value = "Región de O'Higgins"
$ws.Range("S214").Value = 483
$ws.Range("T214").Value = 15

$ws.Range("D215").Value = 44438
$ws.Range("K215").Value = "Fukumoto"
$ws.Range("L215").Value = "Primera"
$ws.Range("M215").Value = 180
$ws.Range("N215").Value = 6000
$ws.Range("O215").Value = 6500
$ws.Range("P215").Value = 6250
$ws.Range("Q215").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R215").Value = "Región de O'Higgins"
$ws.Range("S215").Value = 417
$ws.Range("T215").Value = 15

$ws.Range("D216").Value = 44438
$ws.Range("K216").Value = "Fukumoto"
$ws.Range("L216").Value = "Segunda"
$ws.Range("M216").Value = 120
$ws.Range("N216").Value = 5000
$ws.Range("O216").Value = 5500
$ws.Range("P216").Value = 5250
$ws.Range("Q216").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R216").Value = "Región de O'Higgins"
$ws.Range("S216").Value = 350
$ws.Range("T216").Value = 15

$ws.Range("D217").Value = 44438
$ws.Range("K217").Value = "Navel Late"
$ws.Range("L217").Value = "Primera"
$ws.Range("M217").Value = 180
$ws.Range("N217").Value = 6000
$ws.Range("O217").Value = 6500
$ws.Range("P217").Value = 6250
$ws.Range("Q217").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R217").Value = "Región de O'Higgins"
$ws.Range("S217").Value = 417
$ws.Range("T217").Value = 15

$ws.Range("D218").Value = 44438
$ws.Range("K218").Value = "Navel Late"
$ws.Range("L218").Value = "Segunda"
$ws.Range("M218").Value = 120
$ws.Range("N218").Value = 5000
$ws.Range("O218").Value = 5500
$ws.Range("P218").Value = 5250
$ws.Range("Q218").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R218").Value = "Región de O'Higgins"
$ws.Range("S218").Value = 350
$ws.Range("T218").Value = 15

$ws.Range("D219").Value = 44819
$ws.Range("K219").Value = "Navel Late"
$ws.Range("L219").Value = "Primera"
$ws.Range("M219").Value = 120
$ws.Range("N219").Value = 7000
$ws.Range("O219").Value = 7500
$ws.Range("P219").Value = 7250
$ws.Range("Q219").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R219").Value = "Región de O'Higgins"
$ws.Range("S219").Value = 483
$ws.Range("T219").Value = 15

$ws.Range("D220").Value = 44819
$ws.Range("K220").Value = "Navel Late"
$ws.Range("L220").Value = "Segunda"
$ws.Range("M220").Value = 120
$ws.Range("N220").Value = 6000
$ws.Range("O220").Value = 6500
$ws.Range("P220").Value = 6250
$ws.Range("Q220").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R220").Value = "Región de O'Higgins"
$ws.Range("S220").Value = 417
$ws.Range("T220").Value = 15

$ws.Range("D221").Value = 44397
$ws.Range("K221").Value = "Fukumoto"
$ws.Range("L221").Value = "Primera"
$ws.Range("M221").Value = 120
$ws.Range("N221").Value = 7500
$ws.Range("O221").Value = 8000
$ws.Range("P221").Value = 7750
$ws.Range("Q221").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R221").Value = "Región de O'Higgins"
$ws.Range("S221").Value = 517
$ws.Range("T221").Value = 15

$ws.Range("D222").Value = 44397
$ws.Range("K222").Value = "Fukumoto"
$ws.Range("L222").Value = "Segunda"
$ws.Range("M222").Value = 80
$ws.Range("N222").Value = 6500
$ws.Range("O222").Value = 6500
$ws.Range("P222").Value = 6500
$ws.Range("Q222").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R222").Value = "Región de O'Higgins"
$ws.Range("S222").Value = 433
$ws.Range("T222").Value = 15

$ws.Range("D223").Value = 44323
$ws.Range("K223").Value = "Fukumoto"
$ws.Range("L223").Value = "Primera"
$ws.Range("M223").Value = 120
$ws.Range("N223").Value = 13000
$ws.Range("O223").Value = 14000
$ws.Range("P223").Value = 13500
$ws.Range("Q223").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R223").Value = "Región de O'Higgins"
$ws.Range("S223").Value = 900
$ws.Range("T223").Value = 15

$ws.Range("D224").Value = 44323
$ws.Range("K224").Value = "Fukumoto"
$ws.Range("L224").Value = "Segunda"
$ws.Range("M224").Value = 80
$ws.Range("N224").Value = 12000
$ws.Range("O224").Value = 12000
$ws.Range("P224").Value = 12000
$ws.Range("Q224").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R224").Value = "Región de O'Higgins"
$ws.Range("S224").Value = 800
$ws.Range("T224").Value = 15

$ws.Range("D225").Value = 44474
$ws.Range("K225").Value = "Lane Late"
$ws.Range("L225").Value = "Primera"
$ws.Range("M225").Value = 240
$ws.Range("N225").Value = 6000
$ws.Range("O225").Value = 6500
$ws.Range("P225").Value = 6250
$ws.Range("Q225").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R225").Value = "Región de O'Higgins"
$ws.Range("S225").Value = 417
$ws.Range("T225").Value = 15

$ws.Range("D226").Value = 44474
$ws.Range("K226").Value = "Lane Late"
$ws.Range("L226").Value = "Segunda"
$ws.Range("M226").Value = 200
$ws.Range("N226").Value = 5000
$ws.Range("O226").Value = 5500
$ws.Range("P226").Value = 5250
$ws.Range("Q226").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R226").Value = "Región de O'Higgins"
$ws.Range("S226").Value = 350
$ws.Range("T226").Value = 15

$ws.Range("D227").Value = 44474
$ws.Range("K227").Value = "Navel Late"
$ws.Range("L227").Value = "Primera"
$ws.Range("M227").Value = 240
$ws.Range("N227").Value = 6000
$ws.Range("O227").Value = 6500
$ws.Range("P227").Value = 6250
$ws.Range("Q227").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R227").Value = "Región de O'Higgins"
$ws.Range("S227").Value = 417
$ws.Range("T227").Value = 15

$ws.Range("D228").Value = 44474
$ws.Range("K228").Value = "Navel Late"
$ws.Range("L228").Value = "Segunda"
$ws.Range("M228").Value = 240
$ws.Range("N228").Value = 5000
$ws.Range("O228").Value = 5500
$ws.Range("P228").Value = 5250
$ws.Range("Q228").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R228").Value = "Región de O'Higgins"
$ws.Range("S228").Value = 350
$ws.Range("T228").Value = 15

$ws.Range("D229").Value = 44545
$ws.Range("K229").Value = "Valencia"
$ws.Range("L229").Value = "Primera"
$ws.Range("M229").Value = 60
$ws.Range("N229").Value = 10000
$ws.Range("O229").Value = 11000
$ws.Range("P229").Value = 10500
$ws.Range("Q229").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R229").Value = "Región de O'Higgins"
$ws.Range("S229").Value = 700
$ws.Range("T229").Value = 15

$ws.Range("D230").Value = 44162
$ws.Range("K230").Value = "Valencia"
$ws.Range("L230").Value = "Primera"
$ws.Range("M230").Value = 120
$ws.Range("N230").Value = 14000
$ws.Range("O230").Value = 15000
$ws.Range("P230").Value = 14500
$ws.Range("Q230").Value = "`$/caja 15 kilos granel"
$ws.Range("R230").Value = "Región de O'Higgins"
$ws.Range("S230").Value = 967
$ws.Range("T230").Value = 15

$ws.Range("D231").Value = 44837
$ws.Range("K231").Value = "Navel Late"
$ws.Range("L231").Value = "Primera"
$ws.Range("M231").Value = 120
$ws.Range("N231").Value = 7500
$ws.Range("O231").Value = 8000
$ws.Range("P231").Value = 7750
$ws.Range("Q231").Value = "`$/caja 18 kilos importada"
$ws.Range("R231").Value = "Región de O'Higgins"
$ws.Range("S231").Value = 431
$ws.Range("T231").Value = 18

$ws.Range("D232").Value = 44837
$ws.Range("K232").Value = "Navel Late"
$ws.Range("L232").Value = "Segunda"
$ws.Range("M232").Value = 60
$ws.Range("N232").Value = 6000
$ws.Range("O232").Value = 6000
$ws.Range("P232").Value = 6000
$ws.Range("Q232").Value = "`$/caja 18 kilos importada"
$ws.Range("R232").Value = "Región de O'Higgins"
$ws.Range("S232").Value = 333
$ws.Range("T232").Value = 18

$ws.Range("D233").Value = 44663
$ws.Range("K233").Value = "Valencia"
$ws.Range("L233").Value = "Primera"
$ws.Range("M233").Value = 200
$ws.Range("N233").Value = 10000
$ws.Range("O233").Value = 11000
$ws.Range("P233").Value = 10500
$ws.Range("Q233").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R233").Value = "Región de O'Higgins"
$ws.Range("S233").Value = 700
$ws.Range("T233").Value = 15

$ws.Range("D234").Value = 44704
$ws.Range("K234").Value = "Valencia"
$ws.Range("L234").Value = "Primera"
$ws.Range("M234").Value = 120
$ws.Range("N234").Value = 9000
$ws.Range("O234").Value = 10000
$ws.Range("P234").Value = 9500
$ws.Range("Q234").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R234").Value = "Región de O'Higgins"
$ws.Range("S234").Value = 633
$ws.Range("T234").Value = 15

$ws.Range("D235").Value = 44704
$ws.Range("K235").Value = "Valencia"
$ws.Range("L235").Value = "Segunda"
$ws.Range("M235").Value = 60
$ws.Range("N235").Value = 11000
$ws.Range("O235").Value = 11000
$ws.Range("P235").Value = 11000
$ws.Range("Q235").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R235").Value = "Región de O'Higgins"
$ws.Range("S235").Value = 733
$ws.Range("T235").Value = 15

$ws.Range("D236").Value = 44680
$ws.Range("K236").Value = "Valencia"
$ws.Range("L236").Value = "Primera"
$ws.Range("M236").Value = 120
$ws.Range("N236").Value = 10000
$ws.Range("O236").Value = 11000
$ws.Range("P236").Value = 10500
$ws.Range("Q236").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R236").Value = "Región de O'Higgins"
$ws.Range("S236").Value = 700
$ws.Range("T236").Value = 15

$ws.Range("D237").Value = 44525
$ws.Range("K237").Value = "Navel Late"
$ws.Range("L237").Value = "Primera"
$ws.Range("M237").Value = 120
$ws.Range("N237").Value = 8000
$ws.Range("O237").Value = 9000
$ws.Range("P237").Value = 8500
$ws.Range("Q237").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R237").Value = "Región de O'Higgins"
$ws.Range("S237").Value = 567
$ws.Range("T237").Value = 15

$ws.Range("D238").Value = 44664
$ws.Range("K238").Value = "Valencia"
$ws.Range("L238").Value = "Primera"
$ws.Range("M238").Value = 120
$ws.Range("N238").Value = 10000
$ws.Range("O238").Value = 11000
$ws.Range("P238").Value = 10500
$ws.Range("Q238").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R238").Value = "Región de O'Higgins"
$ws.Range("S238").Value = 700
$ws.Range("T238").Value = 15

$ws.Range("D239").Value = 44218
$ws.Range("K239").Value = "Valencia"
$ws.Range("L239").Value = "Primera"
$ws.Range("M239").Value = 140
$ws.Range("N239").Value = 18000
$ws.Range("O239").Value = 19000
$ws.Range("P239").Value = 18429
$ws.Range("Q239").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R239").Value = "Región de O'Higgins"
$ws.Range("S239").Value = 1229
$ws.Range("T239").Value = 15

$ws.Range("D240").Value = 44421
$ws.Range("K240").Value = "Fukumoto"
$ws.Range("L240").Value = "Primera"
$ws.Range("M240").Value = 200
$ws.Range("N240").Value = 6000
$ws.Range("O240").Value = 6200
$ws.Range("P240").Value = 6100
$ws.Range("Q240").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R240").Value = "Región de O'Higgins"
$ws.Range("S240").Value = 407
$ws.Range("T240").Value = 15

$ws.Range("D241").Value = 44421
$ws.Range("K241").Value = "Fukumoto"
$ws.Range("L241").Value = "Segunda"
$ws.Range("M241").Value = 120
$ws.Range("N241").Value = 5000
$ws.Range("O241").Value = 5500
$ws.Range("P241").Value = 5250
$ws.Range("Q241").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R241").Value = "Región de O'Higgins"
$ws.Range("S241").Value = 350
$ws.Range("T241").Value = 15

$ws.Range("D242").Value = 44421
$ws.Range("K242").Value = "Navel Late"
$ws.Range("L242").Value = "Primera"
$ws.Range("M242").Value = 200
$ws.Range("N242").Value = 5000
$ws.Range("O242").Value = 5500
$ws.Range("P242").Value = 5250
$ws.Range("Q242").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R242").Value = "Región de O'Higgins"
$ws.Range("S242").Value = 350
$ws.Range("T242").Value = 15

$ws.Range("D243").Value = 44421
$ws.Range("K243").Value = "Navel Late"
$ws.Range("L243").Value = "Segunda"
$ws.Range("M243").Value = 160
$ws.Range("N243").Value = 4500
$ws.Range("O243").Value = 4800
$ws.Range("P243").Value = 4650
$ws.Range("Q243").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R243").Value = "Región de O'Higgins"
$ws.Range("S243").Value = 310
$ws.Range("T243").Value = 15

$ws.Range("D244").Value = 44827
$ws.Range("K244").Value = "Navel Late"
$ws.Range("L244").Value = "Primera"
$ws.Range("M244").Value = 120
$ws.Range("N244").Value = 7000
$ws.Range("O244").Value = 7500
$ws.Range("P244").Value = 7250
$ws.Range("Q244").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R244").Value = "Región de O'Higgins"
$ws.Range("S244").Value = 483
$ws.Range("T244").Value = 15

$ws.Range("D245").Value = 44827
$ws.Range("K245").Value = "Navel Late"
$ws.Range("L245").Value = "Segunda"
$ws.Range("M245").Value = 80
$ws.Range("N245").Value = 6000
$ws.Range("O245").Value = 6000
$ws.Range("P245").Value = 6000
$ws.Range("Q245").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R245").Value = "Región de O'Higgins"
$ws.Range("S245").Value = 400
$ws.Range("T245").Value = 15

$ws.Range("D246").Value = 44369
$ws.Range("K246").Value = "Fukumoto"
$ws.Range("L246").Value = "Primera"
$ws.Range("M246").Value = 120
$ws.Range("N246").Value = 8500
$ws.Range("O246").Value = 9000
$ws.Range("P246").Value = 8750
$ws.Range("Q246").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R246").Value = "Región de O'Higgins"
$ws.Range("S246").Value = 583
$ws.Range("T246").Value = 15

$ws.Range("D247").Value = 44589
$ws.Range("K247").Value = "Valencia"
$ws.Range("L247").Value = "Primera"
$ws.Range("M247").Value = 120
$ws.Range("N247").Value = 9000
$ws.Range("O247").Value = 10000
$ws.Range("P247").Value = 9500
$ws.Range("Q247").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R247").Value = "Región de O'Higgins"
$ws.Range("S247").Value = 633
$ws.Range("T247").Value = 15

$ws.Range("D248").Value = 44195
$ws.Range("K248").Value = "Valencia"
$ws.Range("L248").Value = "Primera"
$ws.Range("M248").Value = 140
$ws.Range("N248").Value = 15000
$ws.Range("O248").Value = 16000
$ws.Range("P248").Value = 15571
$ws.Range("Q248").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R248").Value = "Región de O'Higgins"
$ws.Range("S248").Value = 1038
$ws.Range("T248").Value = 15

$ws.Range("D249").Value = 44210
$ws.Range("K249").Value = "Valencia"
$ws.Range("L249").Value = "Primera"
$ws.Range("M249").Value = 115
$ws.Range("N249").Value = 17000
$ws.Range("O249").Value = 18000
$ws.Range("P249").Value = 17609
$ws.Range("Q249").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R249").Value = "Región de O'Higgins"
$ws.Range("S249").Value = 1174
$ws.Range("T249").Value = 15

$ws.Range("D250").Value = 44355
$ws.Range("K250").Value = "Fukumoto"
$ws.Range("L250").Value = "Primera"
$ws.Range("M250").Value = 120
$ws.Range("N250").Value = 8500
$ws.Range("O250").Value = 9000
$ws.Range("P250").Value = 8750
$ws.Range("Q250").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R250").Value = "Región de O'Higgins"
$ws.Range("S250").Value = 583
$ws.Range("T250").Value = 15

$ws.Range("D251").Value = 44159
$ws.Range("K251").Value = "Valencia"
$ws.Range("L251").Value = "Primera"
$ws.Range("M251").Value = 60
$ws.Range("N251").Value = 14000
$ws.Range("O251").Value = 15000
$ws.Range("P251").Value = 14500
$ws.Range("Q251").Value = "`$/caja 15 kilos granel"
$ws.Range("R251").Value = "Región de O'Higgins"
$ws.Range("S251").Value = 967
$ws.Range("T251").Value = 15

$ws.Range("D252").Value = 44795
$ws.Range("K252").Value = "Navel Late"
$ws.Range("L252").Value = "Primera"
$ws.Range("M252").Value = 120
$ws.Range("N252").Value = 5000
$ws.Range("O252").Value = 5500
$ws.Range("P252").Value = 5250
$ws.Range("Q252").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R252").Value = "Región de O'Higgins"
$ws.Range("S252").Value = 350
$ws.Range("T252").Value = 15

$ws.Range("D253").Value = 44777
$ws.Range("K253").Value = "Fukumoto"
$ws.Range("L253").Value = "Primera"
$ws.Range("M253").Value = 120
$ws.Range("N253").Value = 6000
$ws.Range("O253").Value = 6500
$ws.Range("P253").Value = 6250
$ws.Range("Q253").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R253").Value = "Región de O'Higgins"
$ws.Range("S253").Value = 417
$ws.Range("T253").Value = 15

$ws.Range("D254").Value = 44777
$ws.Range("K254").Value = "Fukumoto"
$ws.Range("L254").Value = "Segunda"
$ws.Range("M254").Value = 120
$ws.Range("N254").Value = 5000
$ws.Range("O254").Value = 5500
$ws.Range("P254").Value = 5250
$ws.Range("Q254").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R254").Value = "Región de O'Higgins"
$ws.Range("S254").Value = 350
$ws.Range("T254").Value = 15

$ws.Range("D255").Value = 44453
$ws.Range("K255").Value = "Navel Late"
$ws.Range("L255").Value = "Primera"
$ws.Range("M255").Value = 240
$ws.Range("N255").Value = 6500
$ws.Range("O255").Value = 7000
$ws.Range("P255").Value = 6750
$ws.Range("Q255").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R255").Value = "Región de O'Higgins"
$ws.Range("S255").Value = 450
$ws.Range("T255").Value = 15

$ws.Range("D256").Value = 44453
$ws.Range("K256").Value = "Navel Late"
$ws.Range("L256").Value = "Segunda"
$ws.Range("M256").Value = 180
$ws.Range("N256").Value = 5500
$ws.Range("O256").Value = 6000
$ws.Range("P256").Value = 5750
$ws.Range("Q256").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R256").Value = "Región de O'Higgins"
$ws.Range("S256").Value = 383
$ws.Range("T256").Value = 15

$ws.Range("D257").Value = 44756
$ws.Range("K257").Value = "Fukumoto"
$ws.Range("L257").Value = "Primera"
$ws.Range("M257").Value = 160
$ws.Range("N257").Value = 7500
$ws.Range("O257").Value = 8000
$ws.Range("P257").Value = 7750
$ws.Range("Q257").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R257").Value = "Región de O'Higgins"
$ws.Range("S257").Value = 517
$ws.Range("T257").Value = 15

$ws.Range("D258").Value = 44756
$ws.Range("K258").Value = "Fukumoto"
$ws.Range("L258").Value = "Segunda"
$ws.Range("M258").Value = 80
$ws.Range("N258").Value = 7000
$ws.Range("O258").Value = 7000
$ws.Range("P258").Value = 7000
$ws.Range("Q258").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R258").Value = "Región de O'Higgins"
$ws.Range("S258").Value = 467
$ws.Range("T258").Value = 15

$ws.Range("D259").Value = 44356
$ws.Range("K259").Value = "Fukumoto"
$ws.Range("L259").Value = "Primera"
$ws.Range("M259").Value = 120
$ws.Range("N259").Value = 9000
$ws.Range("O259").Value = 10000
$ws.Range("P259").Value = 9500
$ws.Range("Q259").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R259").Value = "Región de O'Higgins"
$ws.Range("S259").Value = 633
$ws.Range("T259").Value = 15

$ws.Range("D260").Value = 44637
$ws.Range("K260").Value = "Valencia"
$ws.Range("L260").Value = "Primera"
$ws.Range("M260").Value = 120
$ws.Range("N260").Value = 9000
$ws.Range("O260").Value = 10000
$ws.Range("P260").Value = 9500
$ws.Range("Q260").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R260").Value = "Región de O'Higgins"
$ws.Range("S260").Value = 633
$ws.Range("T260").Value = 15

$ws.Range("D261").Value = 44414
$ws.Range("K261").Value = "Fukumoto"
$ws.Range("L261").Value = "Primera"
$ws.Range("M261").Value = 160
$ws.Range("N261").Value = 6000
$ws.Range("O261").Value = 6200
$ws.Range("P261").Value = 6100
$ws.Range("Q261").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R261").Value = "Región de O'Higgins"
$ws.Range("S261").Value = 407
$ws.Range("T261").Value = 15

$ws.Range("D262").Value = 44414
$ws.Range("K262").Value = "Fukumoto"
$ws.Range("L262").Value = "Segunda"
$ws.Range("M262").Value = 120
$ws.Range("N262").Value = 5500
$ws.Range("O262").Value = 5800
$ws.Range("P262").Value = 5650
$ws.Range("Q262").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R262").Value = "Región de O'Higgins"
$ws.Range("S262").Value = 377
$ws.Range("T262").Value = 15

$ws.Range("D263").Value = 44414
$ws.Range("K263").Value = "Navel Late"
$ws.Range("L263").Value = "Primera"
$ws.Range("M263").Value = 300
$ws.Range("N263").Value = 6500
$ws.Range("O263").Value = 7000
$ws.Range("P263").Value = 6750
$ws.Range("Q263").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R263").Value = "Región de O'Higgins"
$ws.Range("S263").Value = 450
$ws.Range("T263").Value = 15

$ws.Range("D264").Value = 44414
$ws.Range("K264").Value = "Navel Late"
$ws.Range("L264").Value = "Segunda"
$ws.Range("M264").Value = 200
$ws.Range("N264").Value = 5500
$ws.Range("O264").Value = 6000
$ws.Range("P264").Value = 5750
$ws.Range("Q264").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R264").Value = "Región de O'Higgins"
$ws.Range("S264").Value = 383
$ws.Range("T264").Value = 15

$ws.Range("D265").Value = 44676
$ws.Range("K265").Value = "Valencia"
$ws.Range("L265").Value = "Primera"
$ws.Range("M265").Value = 120
$ws.Range("N265").Value = 10000
$ws.Range("O265").Value = 11000
$ws.Range("P265").Value = 10500
$ws.Range("Q265").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R265").Value = "Región de O'Higgins"
$ws.Range("S265").Value = 700
$ws.Range("T265").Value = 15

$ws.Range("D266").Value = 44540
$ws.Range("K266").Value = "Navel Late"
$ws.Range("L266").Value = "Primera"
$ws.Range("M266").Value = 200
$ws.Range("N266").Value = 8500
$ws.Range("O266").Value = 9000
$ws.Range("P266").Value = 8750
$ws.Range("Q266").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R266").Value = "Región de O'Higgins"
$ws.Range("S266").Value = 583
$ws.Range("T266").Value = 15

$ws.Range("D267").Value = 44540
$ws.Range("K267").Value = "Valencia"
$ws.Range("L267").Value = "Primera"
$ws.Range("M267").Value = 240
$ws.Range("N267").Value = 8500
$ws.Range("O267").Value = 9000
$ws.Range("P267").Value = 8750
$ws.Range("Q267").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R267").Value = "Región de O'Higgins"
$ws.Range("S267").Value = 583
$ws.Range("T267").Value = 15

$ws.Range("D268").Value = 44383
$ws.Range("K268").Value = "Fukumoto"
$ws.Range("L268").Value = "Primera"
$ws.Range("M268").Value = 120
$ws.Range("N268").Value = 8500
$ws.Range("O268").Value = 9000
$ws.Range("P268").Value = 8750
$ws.Range("Q268").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R268").Value = "Región de O'Higgins"
$ws.Range("S268").Value = 583
$ws.Range("T268").Value = 15

$ws.Range("D269").Value = 44383
$ws.Range("K269").Value = "Fukumoto"
$ws.Range("L269").Value = "Segunda"
$ws.Range("M269").Value = 60
$ws.Range("N269").Value = 7000
$ws.Range("O269").Value = 7000
$ws.Range("P269").Value = 7000
$ws.Range("Q269").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R269").Value = "Región de O'Higgins"
$ws.Range("S269").Value = 467
$ws.Range("T269").Value = 15

$ws.Range("D270").Value = 44692
$ws.Range("K270").Value = "Valencia"
$ws.Range("L270").Value = "Primera"
$ws.Range("M270").Value = 120
$ws.Range("N270").Value = 10000
$ws.Range("O270").Value = 11000
$ws.Range("P270").Value = 10500
$ws.Range("Q270").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R270").Value = "Región de O'Higgins"
$ws.Range("S270").Value = 700
$ws.Range("T270").Value = 15

$ws.Range("D271").Value = 44694
$ws.Range("K271").Value = "Valencia"
$ws.Range("L271").Value = "Primera"
$ws.Range("M271").Value = 120
$ws.Range("N271").Value = 9500
$ws.Range("O271").Value = 10000
$ws.Range("P271").Value = 9750
$ws.Range("Q271").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R271").Value = "Región de O'Higgins"
$ws.Range("S271").Value = 650
$ws.Range("T271").Value = 15

$ws.Range("D272").Value = 44553
$ws.Range("K272").Value = "Valencia"
$ws.Range("L272").Value = "Primera"
$ws.Range("M272").Value = 300
$ws.Range("N272").Value = 9000
$ws.Range("O272").Value = 9500
$ws.Range("P272").Value = 9250
$ws.Range("Q272").Value = "`$/caja 15 kilos granel"
$ws.Range("R272").Value = "Región de O'Higgins"
$ws.Range("S272").Value = 617
$ws.Range("T272").Value = 15

$ws.Range("D273").Value = 44831
$ws.Range("K273").Value = "Navel Late"
$ws.Range("L273").Value = "Primera"
$ws.Range("M273").Value = 120
$ws.Range("N273").Value = 5500
$ws.Range("O273").Value = 6000
$ws.Range("P273").Value = 5750
$ws.Range("Q273").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R273").Value = "Región de O'Higgins"
$ws.Range("S273").Value = 383
$ws.Range("T273").Value = 15

$ws.Range("D274").Value = 44831
$ws.Range("K274").Value = "Navel Late"
$ws.Range("L274").Value = "Segunda"
$ws.Range("M274").Value = 80
$ws.Range("N274").Value = 5000
$ws.Range("O274").Value = 5000
$ws.Range("P274").Value = 5000
$ws.Range("Q274").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R274").Value = "Región de O'Higgins"
$ws.Range("S274").Value = 333
$ws.Range("T274").Value = 15

$ws.Range("D275").Value = 44791
$ws.Range("K275").Value = "Navel Late"
$ws.Range("L275").Value = "Primera"
$ws.Range("M275").Value = 120
$ws.Range("N275").Value = 5500
$ws.Range("O275").Value = 6000
$ws.Range("P275").Value = 5750
$ws.Range("Q275").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R275").Value = "Región de O'Higgins"
$ws.Range("S275").Value = 383
$ws.Range("T275").Value = 15

$ws.Range("D276").Value = 44391
$ws.Range("K276").Value = "Fukumoto"
$ws.Range("L276").Value = "Primera"
$ws.Range("M276").Value = 120
$ws.Range("N276").Value = 8000
$ws.Range("O276").Value = 8500
$ws.Range("P276").Value = 8250
$ws.Range("Q276").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R276").Value = "Región de O'Higgins"
$ws.Range("S276").Value = 550
$ws.Range("T276").Value = 15

$ws.Range("D277").Value = 44391
$ws.Range("K277").Value = "Fukumoto"
$ws.Range("L277").Value = "Segunda"
$ws.Range("M277").Value = 80
$ws.Range("N277").Value = 7000
$ws.Range("O277").Value = 7000
$ws.Range("P277").Value = 7000
$ws.Range("Q277").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R277").Value = "Región de O'Higgins"
$ws.Range("S277").Value = 467
$ws.Range("T277").Value = 15

$ws.Range("D278").Value = 44701
$ws.Range("K278").Value = "Valencia"
$ws.Range("L278").Value = "Primera"
$ws.Range("M278").Value = 120
$ws.Range("N278").Value = 9500
$ws.Range("O278").Value = 10000
$ws.Range("P278").Value = 9750
$ws.Range("Q278").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R278").Value = "Región de O'Higgins"
$ws.Range("S278").Value = 650
$ws.Range("T278").Value = 15

$ws.Range("D279").Value = 44701
$ws.Range("K279").Value = "Valencia"
$ws.Range("L279").Value = "Segunda"
$ws.Range("M279").Value = 80
$ws.Range("N279").Value = 8000
$ws.Range("O279").Value = 8000
$ws.Range("P279").Value = 8000
$ws.Range("Q279").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R279").Value = "Región de O'Higgins"
$ws.Range("S279").Value = 533
$ws.Range("T279").Value = 15

$ws.Range("D280").Value = 44585
$ws.Range("K280").Value = "Valencia"
$ws.Range("L280").Value = "Primera"
$ws.Range("M280").Value = 120
$ws.Range("N280").Value = 9000
$ws.Range("O280").Value = 10000
$ws.Range("P280").Value = 9500
$ws.Range("Q280").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R280").Value = "Región de O'Higgins"
$ws.Range("S280").Value = 633
$ws.Range("T280").Value = 15

$ws.Range("D281").Value = 44585
$ws.Range("K281").Value = "Valencia"
$ws.Range("L281").Value = "Segunda"
$ws.Range("M281").Value = 50
$ws.Range("N281").Value = 8000
$ws.Range("O281").Value = 8000
$ws.Range("P281").Value = 8000
$ws.Range("Q281").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R281").Value = "Región de O'Higgins"
$ws.Range("S281").Value = 533
$ws.Range("T281").Value = 15

$ws.Range("D282").Value = 44160
$ws.Range("K282").Value = "Valencia"
$ws.Range("L282").Value = "Primera"
$ws.Range("M282").Value = 80
$ws.Range("N282").Value = 14000
$ws.Range("O282").Value = 15000
$ws.Range("P282").Value = 14500
$ws.Range("Q282").Value = "`$/caja 15 kilos granel"
$ws.Range("R282").Value = "Región de O'Higgins"
$ws.Range("S282").Value = 967
$ws.Range("T282").Value = 15

$ws.Range("D283").Value = 44826
$ws.Range("K283").Value = "Lane Late"
$ws.Range("L283").Value = "Primera"
$ws.Range("M283").Value = 160
$ws.Range("N283").Value = 6500
$ws.Range("O283").Value = 7000
$ws.Range("P283").Value = 6750
$ws.Range("Q283").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R283").Value = "Región de O'Higgins"
$ws.Range("S283").Value = 450
$ws.Range("T283").Value = 15

$ws.Range("D284").Value = 44376
$ws.Range("K284").Value = "Fukumoto"
$ws.Range("L284").Value = "Primera"
$ws.Range("M284").Value = 120
$ws.Range("N284").Value = 8000
$ws.Range("O284").Value = 8500
$ws.Range("P284").Value = 8250
$ws.Range("Q284").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R284").Value = "Región de O'Higgins"
$ws.Range("S284").Value = 550
$ws.Range("T284").Value = 15

$ws.Range("D285").Value = 44376
$ws.Range("K285").Value = "Fukumoto"
$ws.Range("L285").Value = "Segunda"
$ws.Range("M285").Value = 120
$ws.Range("N285").Value = 7000
$ws.Range("O285").Value = 7500
$ws.Range("P285").Value = 7250
$ws.Range("Q285").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R285").Value = "Región de O'Higgins"
$ws.Range("S285").Value = 483
$ws.Range("T285").Value = 15

$ws.Range("D286").Value = 44784
$ws.Range("K286").Value = "Navel Late"
$ws.Range("L286").Value = "Primera"
$ws.Range("M286").Value = 120
$ws.Range("N286").Value = 6000
$ws.Range("O286").Value = 6500
$ws.Range("P286").Value = 6250
$ws.Range("Q286").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R286").Value = "Provincia de Curicó"
$ws.Range("S286").Value = 417
$ws.Range("T286").Value = 15

$ws.Range("D287").Value = 44784
$ws.Range("K287").Value = "Navel Late"
$ws.Range("L287").Value = "Segunda"
$ws.Range("M287").Value = 60
$ws.Range("N287").Value = 5500
$ws.Range("O287").Value = 5500
$ws.Range("P287").Value = 5500
$ws.Range("Q287").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R287").Value = "Provincia de Curicó"
$ws.Range("S287").Value = 367
$ws.Range("T287").Value = 15

$ws.Range("D288").Value = 44649
$ws.Range("K288").Value = "Valencia"
$ws.Range("L288").Value = "Primera"
$ws.Range("M288").Value = 120
$ws.Range("N288").Value = 10000
$ws.Range("O288").Value = 11000
$ws.Range("P288").Value = 10500
$ws.Range("Q288").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R288").Value = "Región de O'Higgins"
$ws.Range("S288").Value = 700
$ws.Range("T288").Value = 15

$ws.Range("D289").Value = 44649
$ws.Range("K289").Value = "Valencia"
$ws.Range("L289").Value = "Segunda"
$ws.Range("M289").Value = 60
$ws.Range("N289").Value = 9000
$ws.Range("O289").Value = 9000
$ws.Range("P289").Value = 9000
$ws.Range("Q289").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R289").Value = "Región de O'Higgins"
$ws.Range("S289").Value = 600
$ws.Range("T289").Value = 15

$ws.Range("D290").Value = 44165
$ws.Range("K290").Value = "Valencia"
$ws.Range("L290").Value = "Primera"
$ws.Range("M290").Value = 60
$ws.Range("N290").Value = 14000
$ws.Range("O290").Value = 15000
$ws.Range("P290").Value = 14500
$ws.Range("Q290").Value = "`$/caja 15 kilos granel"
$ws.Range("R290").Value = "Región de O'Higgins"
$ws.Range("S290").Value = 967
$ws.Range("T290").Value = 15

$ws.Range("D291").Value = 44419
$ws.Range("K291").Value = "Fukumoto"
$ws.Range("L291").Value = "Primera"
$ws.Range("M291").Value = 120
$ws.Range("N291").Value = 6000
$ws.Range("O291").Value = 6200
$ws.Range("P291").Value = 6100
$ws.Range("Q291").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R291").Value = "Región de O'Higgins"
$ws.Range("S291").Value = 407
$ws.Range("T291").Value = 15

$ws.Range("D292").Value = 44419
$ws.Range("K292").Value = "Fukumoto"
$ws.Range("L292").Value = "Segunda"
$ws.Range("M292").Value = 60
$ws.Range("N292").Value = 5500
$ws.Range("O292").Value = 5800
$ws.Range("P292").Value = 5650
$ws.Range("Q292").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R292").Value = "Región de O'Higgins"
$ws.Range("S292").Value = 377
$ws.Range("T292").Value = 15

$ws.Range("D293").Value = 44419
$ws.Range("K293").Value = "Navel Late"
$ws.Range("L293").Value = "Primera"
$ws.Range("M293").Value = 220
$ws.Range("N293").Value = 5000
$ws.Range("O293").Value = 5500
$ws.Range("P293").Value = 5227
$ws.Range("Q293").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R293").Value = "Región de O'Higgins"
$ws.Range("S293").Value = 348
$ws.Range("T293").Value = 15

$ws.Range("D294").Value = 44419
$ws.Range("K294").Value = "Navel Late"
$ws.Range("L294").Value = "Segunda"
$ws.Range("M294").Value = 110
$ws.Range("N294").Value = 4500
$ws.Range("O294").Value = 4800
$ws.Range("P294").Value = 4664
$ws.Range("Q294").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R294").Value = "Región de O'Higgins"
$ws.Range("S294").Value = 311
$ws.Range("T294").Value = 15

$ws.Range("D295").Value = 44363
$ws.Range("K295").Value = "Fukumoto"
$ws.Range("L295").Value = "Primera"
$ws.Range("M295").Value = 120
$ws.Range("N295").Value = 11000
$ws.Range("O295").Value = 12000
$ws.Range("P295").Value = 11500
$ws.Range("Q295").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R295").Value = "Región de O'Higgins"
$ws.Range("S295").Value = 767
$ws.Range("T295").Value = 15

$ws.Range("D296").Value = 44363
$ws.Range("K296").Value = "Fukumoto"
$ws.Range("L296").Value = "Segunda"
$ws.Range("M296").Value = 60
$ws.Range("N296").Value = 10000
$ws.Range("O296").Value = 10000
$ws.Range("P296").Value = 10000
$ws.Range("Q296").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R296").Value = "Región de O'Higgins"
$ws.Range("S296").Value = 667
$ws.Range("T296").Value = 15

$ws.Range("D297").Value = 44519
$ws.Range("K297").Value = "Navel Late"
$ws.Range("L297").Value = "Primera"
$ws.Range("M297").Value = 160
$ws.Range("N297").Value = 8000
$ws.Range("O297").Value = 9000
$ws.Range("P297").Value = 8500
$ws.Range("Q297").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R297").Value = "Región de O'Higgins"
$ws.Range("S297").Value = 567
$ws.Range("T297").Value = 15

$ws.Range("D484").Value = 44335
$ws.Range("K484").Value = "Fukumoto"
$ws.Range("L484").Value = "Primera"
$ws.Range("M484").Value = 160
$ws.Range("N484").Value = 12000
$ws.Range("O484").Value = 13000
$ws.Range("P484").Value = 12500
$ws.Range("Q484").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R484").Value = "Región de O'Higgins"
$ws.Range("S484").Value = 833
$ws.Range("T484").Value = 15

$ws.Range("D485").Value = 44335
$ws.Range("K485").Value = "Fukumoto"
$ws.Range("L485").Value = "Segunda"
$ws.Range("M485").Value = 60
$ws.Range("N485").Value = 11000
$ws.Range("O485").Value = 11000
$ws.Range("P485").Value = 11000
$ws.Range("Q485").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R485").Value = "Región de O'Higgins"
$ws.Range("S485").Value = 733
$ws.Range("T485").Value = 15

$ws.Rows("486:486").Delete()
